$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133, shifting existing rows 133-225 down to 134-226.
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the "Macroferia Regional de Talca - Coliflor" record.
$ws.Cells.Item(133, 1).Value2 = 5
$ws.Cells.Item(133, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(133, 3).Value2 = "Maule"
$ws.Cells.Item(133, 4).Value2 = 44673
$ws.Cells.Item(133, 5).Value2 = 7
$ws.Cells.Item(133, 6).Value2 = 100112008
$ws.Cells.Item(133, 7).Value2 = "Coliflor"
$ws.Cells.Item(133, 8).Value2 = "Sin especificar"
$ws.Cells.Item(133, 9).Value2 = "Primera"
$ws.Cells.Item(133, 10).Value2 = 2000
$ws.Cells.Item(133, 11).Value2 = 1000
$ws.Cells.Item(133, 12).Value2 = 1000
$ws.Cells.Item(133, 13).Value2 = 1000
$ws.Cells.Item(133, 14).Value2 = "`$/unidad"
$ws.Cells.Item(133, 15).Value2 = "Región del Maule"
$ws.Cells.Item(133, 16).Value2 = 1000
$ws.Cells.Item(133, 17).Value2 = 1
$ws.Cells.Item(133, 18).Value2 = "Hortaliza"
